$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A11").Value = "UI 이미지 적용(퀘스트구슬 표시/ 일시중지 창-다시시작,뒤로가기/ 배경음,효과음켜고끄기)"
$ws.Range("B11").Value = "O"
